# Updates the cached "datetimeFigureOut" field text (11/8/2021 -> 12/8/2021)
# on the slide master and on every slide layout, and updates the
# "Python 3.6+ venv" label on slide 1 to "Python 3.7+ venv".

function Update-DateField {
    param($shapes)
    for ($shpIdx = 1; $shpIdx -le $shapes.Count; $shpIdx++) {
        $shp = $shapes.Item($shpIdx)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "11/8/2021") {
                $tr.Text = "12/8/2021"
            }
        }
    }
}

$p = $ppt.ActivePresentation

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DateField $master.Shapes

# Every slide layout's date placeholder (names/indexes vary per layout).
$layouts = $master.CustomLayouts
Update-DateField $layouts.Item(1).Shapes
Update-DateField $layouts.Item(2).Shapes
Update-DateField $layouts.Item(3).Shapes
Update-DateField $layouts.Item(4).Shapes
Update-DateField $layouts.Item(5).Shapes
Update-DateField $layouts.Item(6).Shapes
Update-DateField $layouts.Item(7).Shapes
Update-DateField $layouts.Item(8).Shapes
Update-DateField $layouts.Item(9).Shapes
Update-DateField $layouts.Item(10).Shapes
Update-DateField $layouts.Item(11).Shapes

# Slide 1: bump the Python version called out in the workflow picture.
$slide = $p.Slides.Item(1)
for ($sIdx = 1; $sIdx -le $slide.Shapes.Count; $sIdx++) {
    $shp = $slide.Shapes.Item($sIdx)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Python 3.6+ venv") {
            $tr.Text = "Python 3.7+ venv"
        }
    }
}
